$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.82
$ws.Range("L2").Value = 1.28
$ws.Range("P2").Value = 2.6
$ws.Range("Q2").Value = 1.56
$ws.Range("R2").Value = 1.68
$ws.Range("S2").Value = 2.34
$ws.Range("T2").Value = 1.52
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 28
$ws.Range("AH2").Value = 14.5
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 29

# Row 3
$ws.Range("P3").Value = 1.29

# Row 4
$ws.Range("H4").Value = 1.87
$ws.Range("I4").Value = 2.42
$ws.Range("J4").Value = 2.46
$ws.Range("K4").Value = 5.8

# Row 5
$ws.Range("F5").Value = 2.78

# Row 7
$ws.Range("F7").Value = 2
$ws.Range("H7").Value = 2.92
$ws.Range("I7").Value = 5.9
$ws.Range("J7").Value = 3.05
$ws.Range("K7").Value = 6.4
$ws.Range("N7").Value = 1.56
$ws.Range("P7").Value = 1.56

# Row 9
$ws.Range("F9").Value = 3.4
$ws.Range("H9").Value = 2.42
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 2.08
$ws.Range("K9").Value = 3.45
$ws.Range("P9").Value = 1.26
$ws.Range("Q9").Value = 3.4

# Row 12
$ws.Range("H12").Value = 2.84
$ws.Range("J12").Value = 3.05
$ws.Range("K12").Value = 5

# Row 13
$ws.Range("J13").Value = 3.7
$ws.Range("P13").Value = 2.08
$ws.Range("Q13").Value = 1.65

# Row 18
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 1.33
$ws.Range("K18").Value = 6.2
$ws.Range("Q18").Value = 1.79

# Row 20/21 swap with updates
$ws.Range("A20").Value = "Chilean Primera Division"
$ws.Range("C20").Value = "18:00:00"
$ws.Range("D20").Value = "Colo Colo"
$ws.Range("E20").Value = "Deportes Limache"
$ws.Range("F20").Value = 1.55
$ws.Range("G20").Value = 980
$ws.Range("H20").Value = 1.04
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 1.03
$ws.Range("K20").Value = 1000
$ws.Range("P20").Value = 1.87
$ws.Range("Q20").Value = 1.74

$ws.Range("A21").Value = "Brazilian Serie B"
$ws.Range("C21").Value = "19:00:00"
$ws.Range("D21").Value = "Chapecoense"
$ws.Range("E21").Value = "Operario PR"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("P21").Value = 1.25
$ws.Range("Q21").Value = 1.01
